$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the relevant paragraphs by their text.
# ------------------------------------------------------------------

$moveParaIndex = 0
$rentParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd("`r`a")
    if ($t -eq "Move player to nearest railroad when card is pulled.") {
        $moveParaIndex = $i
    }
    if ($t -eq "Rent is not being used and asked you to buy it again") {
        $rentParaIndex = $i
    }
}

# ------------------------------------------------------------------
# 1) Remove the sentence "Rent is not being used and asked you to buy
#    it again" (text only; its own paragraph mark is kept), leaving an
#    empty paragraph behind.
# ------------------------------------------------------------------

$rentRange = $d.Paragraphs.Item($rentParaIndex).Range
$rentRange.MoveEnd(1, -1)
$rentRange.Delete()

# ------------------------------------------------------------------
# 2) The paragraph that used to hold the "_GoBack" bookmark is now an
#    empty paragraph right after the (now empty) former "Rent..."
#    paragraph. Delete that trailing empty paragraph mark so the two
#    empty paragraphs collapse into a single one (keeping the former
#    "Rent..." paragraph's own mark/properties, as happens when Word
#    merges a paragraph with the one that follows it).
# ------------------------------------------------------------------

$trailingIndex = $rentParaIndex + 1
$trailingRange = $d.Paragraphs.Item($trailingIndex).Range
$trailingRange.Delete()

# ------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark so that it sits right after the text
#    of the "Move player to nearest railroad when card is pulled."
#    paragraph (it used to live, empty, in the very last paragraph of
#    the document, which was just removed above).
#
#    A collapsed range whose Start/End equals the exact end-of-text
#    offset of a paragraph can't be built directly and reliably (the
#    engine snaps such a range back to the whole paragraph), so we
#    work around it: insert a one-character placeholder right after
#    the sentence, collapse a range immediately in front of that
#    placeholder (which is a perfectly normal, non-boundary position),
#    add the bookmark there, then delete the placeholder character.
# ------------------------------------------------------------------

$moveRange = $d.Paragraphs.Item($moveParaIndex).Range
$moveRange.MoveEnd(1, -1)
$moveRange.InsertAfter("X")

$tempRange = $d.Paragraphs.Item($moveParaIndex).Range
$tempRange.MoveEnd(1, -1)
$tempRange.MoveStart(1, $tempRange.End - $tempRange.Start - 1)
$tempRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $tempRange)

$placeholder = $d.Paragraphs.Item($moveParaIndex).Range
$placeholder.MoveEnd(1, -1)
$placeholder.MoveStart(1, $placeholder.End - $placeholder.Start - 1)
$placeholder.Delete()
